$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1: "This is a Microsoft word document."
#    -> append two trailing spaces (plain) then a parenthetical note in a
#       dark red color (C00000), split across three runs as in the source
#       edit: "(This is a change - Ve" / "rsion for branch alternate" / ")"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$body = $d.Range($p1.Start, $p1.End - 1)

# two plain trailing spaces appended to the existing sentence
$body.InsertAfter("  ")

# the three colored runs that make up the parenthetical comment
$markEnd = $p1.End - 1

$seg1 = "(This is a change " + [char]0x2013 + " Ve"
$ins1 = $d.Range($markEnd, $markEnd)
$ins1.InsertAfter($seg1)
$r1 = $d.Range($markEnd, $markEnd + $seg1.Length)
$r1.Font.Color = 192

$markEnd = $markEnd + $seg1.Length
$seg2 = "rsion for branch alternate"
$ins2 = $d.Range($markEnd, $markEnd)
$ins2.InsertAfter($seg2)
$r2 = $d.Range($markEnd, $markEnd + $seg2.Length)
$r2.Font.Color = 192

$markEnd = $markEnd + $seg2.Length
$seg3 = ")"
$ins3 = $d.Range($markEnd, $markEnd)
$ins3.InsertAfter($seg3)
$r3 = $d.Range($markEnd, $markEnd + $seg3.Length)
$r3.Font.Color = 192

# ---------------------------------------------------------------------------
# 2) Insert a new, empty paragraph right after "It will be treated as a
#    binary file by Git." with a light-grey shading and bold Calibri
#    (paragraph-mark) formatting - matching a blank line copied out of a
#    web page.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$p2.InsertParagraphAfter()

$newPara = $d.Paragraphs(3)
$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 16382457

$markRng = $newPara.Range
$markRng.Font.Name = "Calibri"
$markRng.Font.NameAscii = "Calibri"
$markRng.Font.NameFarEast = "Times New Roman"
$markRng.Font.NameBi = "Calibri"
$markRng.Font.Bold = $true
$markRng.Font.BoldBi = $true
$markRng.Font.Color = 2236704

Write-Output "done"
